# Commit: "Version 1.0! All calculations now functioning correctly. rounds correct!"
#
# The author typed 10 into each of the five "yellow" input cells on the
# Portal sheet (pH, Alkalinity, Calcium, TDS, Temperature). Every other
# cell on Portal/Calculations is a formula and recalculates automatically.

$wb = $excel.ActiveWorkbook

$portal = $wb.Worksheets.Item("Portal")
$portal.Range("C7").Value = 10
$portal.Range("C8").Value = 10
$portal.Range("C9").Value = 10
$portal.Range("C10").Value = 10
$portal.Range("C11").Value = 10

# Match the author's final cursor/viewport position on each sheet.
$calc = $wb.Worksheets.Item("Calculations")
$calc.Activate()
$calc.Range("C46").Select()
$excel.ActiveWindow.Zoom = 85

$portal.Activate()
$portal.Range("C19").Select()
